$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uppercase the surname in B2
$ws.Range("B2").Value = "LESAGE"

# Set the receipt ID in D2
$ws.Range("D2").Value = "250501HL0"
